# Fix Training Data Issue: the BF (Date) column held the malformed text
# "6-6-2012-13" (a leftover artifact of how NBA stats were scraped). It
# should be the correctly formed date string "2013-06-06".
#
# The value must stay literal TEXT (matching the surrounding inlineStr /
# shared-string cells), not get auto-converted into an Excel date serial,
# since "2013-06-06" looks like an ISO date to Excel's input parser.
# Forcing the cell to text format ("@") before the assignment prevents
# that conversion; ClearFormats afterwards drops the temporary explicit
# number format again so the cell ends up with no style override, just
# like the original cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldDate = "6-6-2012-13"
$newDate = "2013-06-06"

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)  # column BF
    if ($cell.Value2 -eq $oldDate) {
        $cell.NumberFormat = "@"
        $cell.Value2 = $newDate
        $cell.ClearFormats()
    }
}
